{"js": "// The document contains one 20-row x 5-col table. The five \"data\" rows\n// (0-indexed: 0, 4, 8, 12, 16) each hold 5 division-practice expressions\n// in their cells; the rows in between are empty spacer rows. We update\n// each cell's expression to its new value in place (scoping the\n// search/replace to that single cell) so formatting (run/paragraph\n// properties) is preserved and so a newly-written value can never be\n// re-matched by a later rule (some new values equal other cells' old\n// values, e.g. \"36\u00f79=\").\nconst rowsData = [\n  // row index, [oldText, newText] for columns 0..4\n  [0, [\"85\u00f77=\", \"62\u00f73=\"], [\"11\u00f78=\", \"12\u00f74=\"], [\"26\u00f79=\", \"77\u00f75=\"], [\"99\u00f78=\", \"78\u00f74=\"], [\"45\u00f78=\", \"79\u00f72=\"]],\n  [4, [\"66\u00f72=\", \"36\u00f79=\"], [\"90\u00f73=\", \"66\u00f78=\"], [\"44\u00f78=\", \"53\u00f79=\"], [\"57\u00f78=\", \"25\u00f73=\"], [\"99\u00f77=\", \"73\u00f77=\"]],\n  [8, [\"98\u00f77=\", \"26\u00f77=\"], [\"19\u00f73=\", \"71\u00f79=\"], [\"10\u00f73=\", \"44\u00f74=\"], [\"36\u00f79=\", \"11\u00f77=\"], [\"17\u00f76=\", \"82\u00f72=\"]],\n  [12, [\"19\u00f74=\", \"13\u00f79=\"], [\"74\u00f75=\", \"27\u00f79=\"], [\"34\u00f78=\", \"41\u00f77=\"], [\"90\u00f77=\", \"35\u00f78=\"], [\"84\u00f78=\", \"97\u00f74=\"]],\n  [16, [\"59\u00f74=\", \"26\u00f76=\"], [\"62\u00f78=\", \"66\u00f79=\"], [\"48\u00f72=\", \"72\u00f76=\"], [\"56\u00f78=\", \"27\u00f75=\"], [\"39\u00f79=\", \"21\u00f79=\"]],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (const [rowIndex, ...cells] of rowsData) {\n  for (let col = 0; col < cells.length; col++) {\n    const [oldText, newText] = cells[col];\n    const cell = table.getCell(rowIndex, col);\n    const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (const range of results.items) {\n      range.insertText(newText, Word.InsertLocation.replace);\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# The document contains one 20-row x 5-col table. The five \"data\" rows\n# (1-indexed table rows: 1, 5, 9, 13, 17) each hold 5 division-practice\n# expressions in their cells; the rows in between are empty spacer rows.\n# We set each cell's Range.Text directly (row/column addressed, not a\n# text search) so every cell is targeted unambiguously even though some\n# new values equal other cells' old values (e.g. \"36\u00f79=\").\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowIndexes = @(1, 5, 9, 13, 17)\n$rowValues = @(\n  @(\"62\u00f73=\", \"12\u00f74=\", \"77\u00f75=\", \"78\u00f74=\", \"79\u00f72=\"),\n  @(\"36\u00f79=\", \"66\u00f78=\", \"53\u00f79=\", \"25\u00f73=\", \"73\u00f77=\"),\n  @(\"26\u00f77=\", \"71\u00f79=\", \"44\u00f74=\", \"11\u00f77=\", \"82\u00f72=\"),\n  @(\"13\u00f79=\", \"27\u00f79=\", \"41\u00f77=\", \"35\u00f78=\", \"97\u00f74=\"),\n  @(\"26\u00f76=\", \"66\u00f79=\", \"72\u00f76=\", \"27\u00f75=\", \"21\u00f79=\")\n)\n\nfor ($i = 0; $i -lt $rowIndexes.Count; $i++) {\n  $row = $rowIndexes[$i]\n  $values = $rowValues[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $t.Cell($row, $col).Range.Text = $values[$col - 1]\n  }\n}\n"}
